$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update existing D column values (rows 2-9)
$ws.Range("D2").Value = 1.005
$ws.Range("D3").Value = 1.005
$ws.Range("D4").Value = 171.513140029589
$ws.Range("D5").Value = 273.818043956044
$ws.Range("D6").Value = 170.39500899668
$ws.Range("D7").Value = 35.6017473819737
$ws.Range("D8").Value = 2.20316449136851
$ws.Range("D9").Value = 0.0215224998099757

# Add new row 10
$ws.Range("A10").Value = 7
$ws.Range("B10").Value = 662.552639418119
$ws.Range("C10").Value = -0.00000000000500222085975111
$ws.Range("D10").Value = 0.0000120626532407186
